$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values that live on row 39 for the
# zib-side columns (B,C,D,E,G,H,I,J) so they can be moved down to row 41.
# NOTE: ".Value" getter is unreliable in this runtime, so ".Value2" is used
# for reads; both work fine for writes.
$b39 = $ws.Range("B39").Value2
$c39 = $ws.Range("C39").Value2
$d39 = $ws.Range("D39").Value2
$e39 = $ws.Range("E39").Value2
$g39 = $ws.Range("G39").Value2
$h39 = $ws.Range("H39").Value2
$i39 = $ws.Range("I39").Value2
$j39 = $ws.Range("J39").Value2

# Clear row 39's zib-side columns (B,C,D,E,G,H,I,J) -- these become blank.
$ws.Range("B39").Value2 = ""
$ws.Range("C39").Value2 = ""
$ws.Range("D39").Value2 = ""
$ws.Range("E39").Value2 = ""
$ws.Range("G39").Value2 = ""
$ws.Range("H39").Value2 = ""
$ws.Range("I39").Value2 = ""
$ws.Range("J39").Value2 = ""

# Move the captured values down onto row 41's zib-side columns.
# Column G holds a cardinality value that looks numeric ("1") but must stay
# text, so force a text number format before assigning it.
$ws.Range("B41").Value2 = $b39
$ws.Range("C41").Value2 = $c39
$ws.Range("D41").Value2 = $d39
$ws.Range("E41").Value2 = $e39
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value2 = $g39
$ws.Range("G41").Style = "Normal"
$ws.Range("H41").Value2 = $h39
$ws.Range("I41").Value2 = $i39
$ws.Range("J41").Value2 = $j39
